$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.098.43'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '2.285.78'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'533.76"
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("D6").Value = "'130.95"
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = '  +4.07%  '
$ws.Range("D9").Value = '2.286.12'
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("D10").Value = "'0.0997"
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").Value = "'0.332"
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").Value = "'23.39"
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("D15").Value = '2.693.10'
$ws.Range("D16").Value = '58.031.70'
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = '2.296.17'
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").Value = "'10.47"
$ws.Range("E19").Value = '  -1.29%  '
$ws.Range("E20").Value = '  -2.27%  '
$ws.Range("D21").Value = "'313.16"
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").Value = "'62.89"
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("E28").Value = '  -3.42%  '
$ws.Range("D29").Value = "'170.62"
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = "'1.06"
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").Value = "'0.377"
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = "'17.81"
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("E40").Value = '  -1.42%  '
$ws.Range("D41").Value = "'140.01"
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("D42").Value = "'286.35"
$ws.Range("E42").Value = '  -2.83%  '
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = "'0.0953"
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("E45").Value = '  -0.52%  '
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = "'17.99"
$ws.Range("E47").Value = '  -0.90%  '
$ws.Range("E48").Value = '  -1.45%  '
$ws.Range("D49").Value = "'10.94"
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("E51").Value = '  +0.61%  '
